# "unify functions into single_run.m + finish roy's review"
#
# The auditory-motor run sheet (Sheet1) is rewritten:
#  - rows 2-5 (run 1, blocks 1-4) get placeholder start_time/play_duration
#    values of 999 (pending re-export from the unified single_run.m script)
#  - row 4's "hand" value is corrected from L -> R (Roy's review)
#  - rows 6-9 (run 1, blocks 5-8) are blanked out to 0 / empty pending rerun
#  - rows 10-33 are left exactly as they were
#  - column C (start_time) narrows now that it just holds "999" placeholders

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row (text unchanged)
$ws.Range("A1").Value = "run_num"
$ws.Range("B1").Value = "block_num"
$ws.Range("C1").Value = "start_time"
$ws.Range("D1").Value = "play_duration"
$ws.Range("E1").Value = "ear"
$ws.Range("F1").Value = "hand"

# Row 2 (run 1, block 1)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 999
$ws.Range("D2").Value = 999
$ws.Range("E2").Value = "R"
$ws.Range("F2").Value = "R"

# Row 3 (run 1, block 2)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 999
$ws.Range("D3").Value = 999
$ws.Range("E3").Value = "R"
$ws.Range("F3").Value = "L"

# Row 4 (run 1, block 3) - hand corrected L -> R per Roy's review
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 999
$ws.Range("D4").Value = 999
$ws.Range("E4").Value = "R"
$ws.Range("F4").Value = "R"

# Row 5 (run 1, block 4)
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 999
$ws.Range("D5").Value = 999
$ws.Range("E5").Value = "R"
$ws.Range("F5").Value = "L"

# Rows 6-9 (run 1, blocks 5-8) blanked out to match the not-yet-rerun rows below
foreach ($r in 6..9) {
    $ws.Range("A$r").Value = 0
    $ws.Range("B$r").Value = 0
    $ws.Range("C$r").Value = 0
    $ws.Range("D$r").Value = 0
    $ws.Range("E$r").ClearContents()
    $ws.Range("F$r").ClearContents()
}

# Column C (start_time) re-fits to its new, shorter "999" placeholder content
$ws.Columns.Item(3).ColumnWidth = 9.5
